# Update "想去人数" (F column) values across sheets, per the source update.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 1048
$ws1.Range("F7").Value  = 10825
$ws1.Range("F10").Value = 301
$ws1.Range("F11").Value = 1051
$ws1.Range("F12").Value = 728
$ws1.Range("F13").Value = 12140
$ws1.Range("F14").Value = 12605
$ws1.Range("F18").Value = 28

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 12
$ws4.Range("F7").Value  = 1048
$ws4.Range("F8").Value  = 10825
$ws4.Range("F11").Value = 301
$ws4.Range("F12").Value = 1051
$ws4.Range("F13").Value = 728
$ws4.Range("F14").Value = 12140
$ws4.Range("F15").Value = 12605
$ws4.Range("F19").Value = 28
